# Swap the data of rows 2 and 3 (all cells that differ between them) on the
# active worksheet. Columns that are identical between the two rows are left
# untouched so their original representation (e.g. text dates) is preserved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values differ between row 2 and row 3 and therefore need to
# be swapped.
$columns = @("A","B","D","E","F","G","H","K","Q","R","AH","AJ","AK","AO")

foreach ($col in $columns) {
    $cell2 = $ws.Range("$col" + "2")
    $cell3 = $ws.Range("$col" + "3")

    $val2 = $cell2.Value2
    $val3 = $cell3.Value2

    $cell2.Value2 = $val3
    $cell3.Value2 = $val2
}
